$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1235.4667
$ws.Range("J17").Value = 1303
$ws.Range("L17").Value = 3909
$ws.Range("N17").Value = -4245

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = 0

$ws.Range("H70").Value = 3537.2856
$ws.Range("I70").Value = 2195.1667
$ws.Range("K70").Value = 6585.500100000001
$ws.Range("M70").Value = -6315.500100000001

$ws.Range("H73").Value = 3537.2856
$ws.Range("I73").Value = 2195.1667
$ws.Range("K73").Value = 6585.500100000001
$ws.Range("M73").Value = -5649.500100000001

$ws.Range("H76").Value = 9099254
$ws.Range("J76").Value = 10790
$ws.Range("L76").Value = 10790
$ws.Range("N76").Value = -11420

$ws.Range("H79").Value = 9099254
$ws.Range("J79").Value = 10790
$ws.Range("L79").Value = 10790
$ws.Range("N79").Value = -12974

$ws.Range("H132").Value = 1556.5264
$ws.Range("I132").Value = 1504.4706
$ws.Range("K132").Value = 4513.4118
$ws.Range("M132").Value = -1983.4118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23343.422
$ws.Range("I32").Value = 24880.143
$ws.Range("J32").Value = 1829.3334
$ws.Range("K32").Value = 24880.143
$ws.Range("L32").Value = 1829.3334
$ws.Range("M32").Value = -24593.143
$ws.Range("N32").Value = -2403.3334

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0

$ws.Range("H61").Value = 6796.3335
$ws.Range("I61").Value = 4231.533
$ws.Range("K61").Value = 4231.533
$ws.Range("M61").Value = -4019.533

$ws.Range("H107").Value = 99969
$ws.Range("J107").Value = 99969
$ws.Range("L107").Value = 99969
$ws.Range("N107").Value = -107649

$ws.Range("H122").Value = 3632.5
$ws.Range("I122").Value = 3527.3076
$ws.Range("K122").Value = 10581.9228
$ws.Range("M122").Value = -8131.9228

$ws.Range("H132").Value = 5456.3105
$ws.Range("I132").Value = 3788.4167
$ws.Range("K132").Value = 11365.2501
$ws.Range("M132").Value = -8835.250100000001

$ws.Range("H136").Value = 6796.3335
$ws.Range("I136").Value = 4231.533
$ws.Range("K136").Value = 12694.599
$ws.Range("M136").Value = -10144.599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 682.6667
$ws.Range("I12").Value = 89.8
$ws.Range("J12").Value = 1423.75
$ws.Range("K12").Value = 89.8
$ws.Range("L12").Value = 1423.75
$ws.Range("M12").Value = 78.2
$ws.Range("N12").Value = -1759.75

$ws.Range("H35").Value = 55000
$ws.Range("J35").Value = 60000
$ws.Range("L35").Value = 60000
$ws.Range("N35").Value = -60620

$ws.Range("H46").Value = 24057
$ws.Range("J46").Value = 24057
$ws.Range("L46").Value = 24057
$ws.Range("N46").Value = -24653

$ws.Range("H56").Value = 17000
$ws.Range("J56").Value = 17000
$ws.Range("L56").Value = 17000
$ws.Range("N56").Value = -18478

$ws.Range("H134").Value = 5423.8667
$ws.Range("I134").Value = 2284.9443
$ws.Range("K134").Value = 6854.8329
$ws.Range("M134").Value = -4319.8329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 195
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 90
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 90
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -370

$ws.Range("H16").Value = 6451.143
$ws.Range("I16").Value = 6315.25
$ws.Range("J16").Value = 6632.3335
$ws.Range("K16").Value = 6315.25
$ws.Range("L16").Value = 6632.3335
$ws.Range("M16").Value = -6028.25
$ws.Range("N16").Value = -7206.3335

$ws.Range("H31").Value = 35718684
$ws.Range("I31").Value = 100001030
$ws.Range("K31").Value = 100001030
$ws.Range("M31").Value = -100000735

$ws.Range("H34").Value = 35718684
$ws.Range("I34").Value = 100001030
$ws.Range("K34").Value = 100001030
$ws.Range("M34").Value = -100000828

$ws.Range("H113").Value = 6451.143
$ws.Range("I113").Value = 6315.25
$ws.Range("J113").Value = 6632.3335
$ws.Range("K113").Value = 6315.25
$ws.Range("L113").Value = 6632.3335
$ws.Range("M113").Value = -4145.25
$ws.Range("N113").Value = -10972.3335

$ws.Range("H134").Value = 11257.934
$ws.Range("J134").Value = 11628.25
$ws.Range("L134").Value = 34884.75
$ws.Range("N134").Value = -39954.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3778.5715
$ws.Range("I80").Value = 1950
$ws.Range("K80").Value = 5850
$ws.Range("M80").Value = -4914

$ws.Range("H83").Value = 3778.5715
$ws.Range("I83").Value = 1950
$ws.Range("K83").Value = 17550
$ws.Range("M83").Value = -12870

$ws.Range("H114").Value = 341
$ws.Range("I114").Value = 78.75
$ws.Range("J114").Value = 865.5
$ws.Range("K114").Value = 236.25
$ws.Range("L114").Value = 2596.5
$ws.Range("M114").Value = 3017.75
$ws.Range("N114").Value = -9104.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 59999
$ws.Range("J32").Value = 59999
$ws.Range("L32").Value = 59999
$ws.Range("N32").Value = -60591

$ws.Range("H42").Value = 65290
$ws.Range("J42").Value = 65290
$ws.Range("L42").Value = 65290
$ws.Range("N42").Value = -66260

$ws.Range("H57").Value = 35965.145
$ws.Range("I57").Value = 22666.666
$ws.Range("J57").Value = 45939
$ws.Range("K57").Value = 22666.666
$ws.Range("L57").Value = 45939
$ws.Range("M57").Value = -21846.666
$ws.Range("N57").Value = -47579

$ws.Range("H70").Value = 18839.24
$ws.Range("I70").Value = 15227.723
$ws.Range("J70").Value = 28126
$ws.Range("K70").Value = 15227.723
$ws.Range("L70").Value = 28126
$ws.Range("M70").Value = -14957.723
$ws.Range("N70").Value = -28666

$ws.Range("H73").Value = 18839.24
$ws.Range("I73").Value = 15227.723
$ws.Range("J73").Value = 28126
$ws.Range("K73").Value = 15227.723
$ws.Range("L73").Value = 28126
$ws.Range("M73").Value = -14291.723
$ws.Range("N73").Value = -29998

$ws.Range("H80").Value = 4144.25
$ws.Range("J80").Value = 4677.5713
$ws.Range("L80").Value = 4677.5713
$ws.Range("N80").Value = -6673.5713

$ws.Range("H83").Value = 4144.25
$ws.Range("J83").Value = 4677.5713
$ws.Range("L83").Value = 23387.8565
$ws.Range("N83").Value = -33371.85649999999

$ws.Range("H115").Value = 65290
$ws.Range("J115").Value = 65290
$ws.Range("L115").Value = 65290
$ws.Range("N115").Value = -67640

$ws.Range("H126").Value = 3020.037
$ws.Range("I126").Value = 2616.923
$ws.Range("K126").Value = 7850.768999999999
$ws.Range("M126").Value = -5380.768999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 24752
$ws.Range("J3").Value = 24752
$ws.Range("L3").Value = 24752
$ws.Range("N3").Value = -24976

$ws.Range("H15").Value = 24752
$ws.Range("J15").Value = 24752
$ws.Range("L15").Value = 24752
$ws.Range("N15").Value = -25092

$ws.Range("H16").Value = 1785.1111
$ws.Range("I16").Value = 1258.25
$ws.Range("K16").Value = 1258.25
$ws.Range("M16").Value = -1088.25

$ws.Range("H46").Value = 3590.1702
$ws.Range("J46").Value = 4162.0264
$ws.Range("L46").Value = 4162.0264
$ws.Range("N46").Value = -4538.0264

$ws.Range("H100").Value = 22729598
$ws.Range("I100").Value = 125000750
$ws.Range("J100").Value = 2676.3333
$ws.Range("K100").Value = 125000750
$ws.Range("L100").Value = 2676.3333
$ws.Range("M100").Value = -125000209
$ws.Range("N100").Value = -3758.3333

$ws.Range("H122").Value = 4116.9546
$ws.Range("I122").Value = 3920.7222
$ws.Range("K122").Value = 11762.1666
$ws.Range("M122").Value = -9312.1666

$ws.Range("H136").Value = 4517.892
$ws.Range("I136").Value = 3084.8572
$ws.Range("J136").Value = 6398.75
$ws.Range("K136").Value = 9254.571599999999
$ws.Range("L136").Value = 19196.25
$ws.Range("M136").Value = -6704.571599999999
$ws.Range("N136").Value = -24296.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 2000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0

$ws.Range("H38").Value = 20995
$ws.Range("I38").Value = 25660
$ws.Range("J38").Value = 7000
$ws.Range("K38").Value = 25660
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = -25187
$ws.Range("N38").Value = -7946

$ws.Range("H55").Value = 1458.75
$ws.Range("J55").Value = 1458.75
$ws.Range("L55").Value = 1458.75
$ws.Range("N55").Value = -2012.75

$ws.Range("H100").Value = 4699.5
$ws.Range("I100").Value = 5481.1816
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 10962.3632
$ws.Range("L100").Value = 3666.6668
$ws.Range("M100").Value = -10421.3632
$ws.Range("N100").Value = -4748.6668

$ws.Range("H115").Value = 76999.60000000001
$ws.Range("J115").Value = 76999.60000000001
$ws.Range("L115").Value = 76999.60000000001
$ws.Range("N115").Value = -80133.60000000001

Write-Output "applied changes"